$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2026-01-24 18:27:25"

# Insert two new blank rows at row 4, pushing the old rows 4-5 down to 6-7.
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# Refresh the timestamp on every data row (it reflects the moment of this scrape).
$ws.Range("A2").Value = $timestamp
$ws.Range("A3").Value = $timestamp
$ws.Range("A4").Value = $timestamp
$ws.Range("A5").Value = $timestamp
$ws.Range("A6").Value = $timestamp
$ws.Range("A7").Value = $timestamp

# New row 4: newly scraped listing.
$ws.Range("B4").Value = "【初心者・未経験OK】 AIを学びながら在宅で働くお仕事|月15〜20万円可能|スマホOK"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5478263"
$ws.Range("G4").Value = 303
$ws.Range("H4").Value = "🔥AI,Ai"

# New row 5: newly scraped listing.
$ws.Range("B5").Value = "「飲み会調整・店舗共有・終電管理・近距離マッチングを備えた飲み会支援アプリの開発依頼」"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5478300"
$ws.Range("G5").Value = 135
$ws.Range("H5").Value = "◆開発 ◇アプリ"

# (The row inserts above already copied the hyperlink cell style ("s=1")
# from the neighbouring F column cells onto F4/F5, matching F2/F3/F6/F7.)

# Rebuild the hyperlinks collection: drop the stale entries (they did not
# follow the row insert) and re-add one per data row in order.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5477958") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5477903") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5478263") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5478300") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5477985") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5418064") | Out-Null

Write-Host "edit complete"
